$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extent on the active sheet.
$lastRow = $ws.Cells.Item(1,1).Worksheet.UsedRange.Rows.Count
$lastCol = $ws.Cells.Item(1,1).Worksheet.UsedRange.Columns.Count

# --- 1) "Förändrad" (column C) timestamp bump: 45184 -> 45186 for every data row ---
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# --- 2) Add the visible link text (the "Beteckning" in column A) as the second
#        HYPERLINK() argument for every hyperlink formula on the sheet (columns
#        S..Y hold: Artfyndslänk, Kartlänk, Knärotsbufferlänk, Klagomålslänk,
#        Klagomålsmaillänk, Tillsynsbegäranslänk, Tillsynsmaillänk). ---
$cols = @("S","T","U","V","W","X","Y")
for ($row = 2; $row -le $lastRow; $row++) {
    $id = $ws.Range("A" + $row).Value2
    foreach ($col in $cols) {
        $cell = $ws.Range($col + $row)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f.IndexOf("HYPERLINK(") -ge 0 -and $f.IndexOf(",") -lt 0) {
                $trimmed = $f.Substring(0, $f.Length - 1)
                $cell.Formula = $trimmed + ', "' + $id + '")'
            }
        }
    }
}
